# Input data change: OperationScenario_Component_HeatingElement
# Row 2 (ID_HeatingElement = 1) -> "power" column value changed from 100000 to 0.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 0
